$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.738.45'
$ws.Range('E2').Value = '  -6.73%  '
$ws.Range('D3').Value = '2.594.55'
$ws.Range('E3').Value = '  -1.28%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.83'
$ws.Range('E5').Value = '  -2.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '96.63'
$ws.Range('E6').Value = '  -5.07%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.583'
$ws.Range('E7').Value = '  -3.80%  '
$ws.Range('E8').Value = '  +0.06%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.562'
$ws.Range('E9').Value = '  -3.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.76'
$ws.Range('E10').Value = '  -8.24%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0820'
$ws.Range('E11').Value = '  -4.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.85'
$ws.Range('E12').Value = '  -5.26%  '
$ws.Range('D13').Value = '2.986.98'
$ws.Range('E13').Value = '  -1.32%  '
$ws.Range('E14').Value = '  +1.05%  '
$ws.Range('D15').Value = '2.582.31'
$ws.Range('E15').Value = '  -1.84%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.898'
$ws.Range('E16').Value = '  -4.05%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.40'
$ws.Range('E17').Value = '  -4.65%  '
$ws.Range('D18').Value = '43.704.11'
$ws.Range('E18').Value = '  -7.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.75'
$ws.Range('E19').Value = '  -0.59%  '
$ws.Range('D20').Value = '0.0₃0985'
$ws.Range('E20').Value = '  -3.82%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.46'
$ws.Range('E21').Value = '  -5.09%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '73.64'
$ws.Range('E22').Value = '  +2.16%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '266.29'
$ws.Range('E23').Value = '  -4.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.95'
$ws.Range('E24').Value = '  -3.69%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.22'
$ws.Range('E25').Value = '  +1.08%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '29.29'
$ws.Range('E26').Value = '  -0.44%  '
$ws.Range('E27').Value = '  +0.06%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.33'
$ws.Range('E28').Value = '  -3.85%  '
$ws.Range('B29').Value = 'InjectiveProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '38.23'
$ws.Range('E29').Value = '  -3.80%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.16'
$ws.Range('E30').Value = '  -6.65%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.16'
$ws.Range('E31').Value = '  -4.43%  '
$ws.Range('E32').Value = '  -1.59%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.23'
$ws.Range('E33').Value = '  -1.77%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '152.30'
$ws.Range('E34').Value = '  -0.27%  '
$ws.Range('E35').Value = '  -3.62%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0823'
$ws.Range('E36').Value = '  -3.29%  '
$ws.Range('E37').Value = '  -5.87%  '
$ws.Range('E38').Value = '  -2.19%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '24.24'
$ws.Range('E39').Value = '  +3.59%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '17.07'
$ws.Range('E40').Value = '  +5.54%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.64'
$ws.Range('E41').Value = '  -2.23%  '
$ws.Range('E42').Value = '  -5.69%  '
$ws.Range('E43').Value = '  -6.09%  '
$ws.Range('D44').Value = '2.040.59'
$ws.Range('E44').Value = '  -4.34%  '
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.31'
$ws.Range('E46').Value = '  -6.29%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.16'
$ws.Range('E47').Value = '  -4.68%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.62'
$ws.Range('E48').Value = '  +4.87%  '
$ws.Range('D49').Value = '2.841.41'
$ws.Range('E49').Value = '  -1.50%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '106.14'
$ws.Range('E50').Value = '  -3.69%  '
$ws.Range('E51').Value = '  -5.35%  '
